$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that became blank (naive forecaster bug fix removed stray values)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Update recomputed forecast values (floating point refinement)
$ws.Range("E3").Value = 10.06916370210016
$ws.Range("C5").Value = -20.40984652067477
$ws.Range("E5").Value = -33.45158075171639
$ws.Range("C6").Value = -14.45332333832744
$ws.Range("E6").Value = 7.857938327064207
$ws.Range("C7").Value = 5.331710924091815
$ws.Range("C8").Value = 8.600536527919612
$ws.Range("C9").Value = 11.04982736891555
$ws.Range("E9").Value = 9.5899211611429
$ws.Range("E11").Value = 7.086193663491014
$ws.Range("C12").Value = 4.639893381363192
$ws.Range("C13").Value = -2.313034291448768
$ws.Range("C17").Value = 3.860244074450203
$ws.Range("E17").Value = 3.254220449867051
$ws.Range("E18").Value = 1.985659800779915
$ws.Range("C21").Value = 4.421855465610269
$ws.Range("E21").Value = 5.26036486209962
$ws.Range("C22").Value = 4.695933104194361
$ws.Range("E22").Value = 6.493919935864612
$ws.Range("E23").Value = 11.40563592910599
$ws.Range("C24").Value = 3.320585727896574
$ws.Range("E24").Value = -4.09821347263859
$ws.Range("C25").Value = 5.276665321936491
$ws.Range("E25").Value = 2.513980370944258
$ws.Range("E27").Value = 2.636028935395318
$ws.Range("C28").Value = 1.782333336406383
$ws.Range("C29").Value = 0.3252781783188441
$ws.Range("E30").Value = 3.191985284262278
$ws.Range("C31").Value = -0.1289008616491394
$ws.Range("E31").Value = -0.9308772335758664
$ws.Range("E32").Value = -11.8352240479
$ws.Range("E34").Value = 53.94004854052483
$ws.Range("C36").Value = 6.928818429977723
$ws.Range("C37").Value = 5.555562147330284
$ws.Range("E39").Value = 20.55431402884296
$ws.Range("C41").Value = 4.273664635170782
$ws.Range("C42").Value = 5.120680133083622
$ws.Range("C43").Value = -0.7317253369667154
$ws.Range("E43").Value = -3.774469028318805
$ws.Range("C44").Value = 0.8562564928550342
$ws.Range("E44").Value = 1.697198638953612
$ws.Range("E45").Value = -3.897546227660653
$ws.Range("C46").Value = -0.5532735011319123
$ws.Range("E46").Value = -2.911323063974536
$ws.Range("C47").Value = -4.911273445035658
$ws.Range("E47").Value = -6.226959158104972
$ws.Range("C48").Value = -0.1644433828108638
$ws.Range("E49").Value = -0.2465826728646148
$ws.Range("E50").Value = -5.499724587330512
$ws.Range("E51").Value = -8.666940084126084
$ws.Range("E53").Value = 0.1329973408891627
